$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 548
    $ws.Range("F4").Value = 273
    $ws.Range("F6").Value = 89
    $ws.Range("F7").Value = 786
}
